$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7: "Experimental" value was blank -> "false"
# A plain .Value = "false" assignment gets auto-coerced to a Boolean by
# Excel's type inference (same as typing FALSE into a cell), so instead
# build it as a text formula and flatten it to a literal value via
# copy / paste-special-values, which preserves the plain string type.
$ws.Range("B7").Formula = "=T(""false"")"
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# Row 8: "Date" value updated
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# Row 17: "Description" value was blank -> new text
$ws.Range("B17").Value = "Trends in VO2max values over time"
